$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.861.96'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.632.55'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.38%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.26'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5098'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2551'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06333'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.37'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.23%  '
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.646.34'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.252'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5405'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0₅7707'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.93'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.898.84'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.004'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '194.44'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.26%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.397'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.890'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.005'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.005'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.854'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.11'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1203'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +5.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.803'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.55'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.88%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.233'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.04899'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.222'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.138'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.524'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.370'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.8856'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.572'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.128.81'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5369'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01545'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.76%  '
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.543'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8104'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.470'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.03'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.81%  '
$ws.Range("E45").Value = '  +2.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.772.87'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4531'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.005'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.44'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("E50").Value = '  +0.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.005'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.08%  '
